# delete visual_clutter and restructure stimulus folder
#
# Row 1 (headers) is unchanged: posFile | stimFile | cs_plus_s | cs_minus_s | cs_plus_ns | cs_minus_ns
# Row 2 (values)  keeps posFile/stimFile entries, but the four stimulus-image paths are
# replaced with the new, flattened "stimuli/<Name>.png" filenames.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the stimulus file paths on row 2 (columns C:F) -----------------
$ws.Range("C2").Value = "stimuli/Social1.png"
$ws.Range("D2").Value = "stimuli/Social2.png"
$ws.Range("E2").Value = "stimuli/Nonsocial2.png"
$ws.Range("F2").Value = "stimuli/Nonsocial1.png"

# --- Resize columns: C:D now share A/B's width, E:F get a bit wider --------
$ws.Columns.Item(3).ColumnWidth = 15.67
$ws.Columns.Item(4).ColumnWidth = 15.67
$ws.Columns.Item(5).ColumnWidth = 18.92
$ws.Columns.Item(6).ColumnWidth = 18.92

# --- Move the active selection to E6 (single cell, nothing selected range) -
[void]$ws.Range("E6").Select()
